# "add basic combo system"
#
# The underlying sheet content/data is unchanged in this revision; what
# changed is purely presentational: several columns on the "List of items"
# sheet were widened (Format > AutoFit Column Width after some cells picked
# up longer text / a new Deity column came into use) and the user's last
# selection moved to H23.
#
# NOTE on the ColumnWidth numbers below: Excel's ColumnWidth property is
# expressed in "characters of the Normal style font" and gets snapped to a
# pixel grid on write, so the literal target widths (which were captured
# from a real Excel session, in fractional character units derived from
# true font-metrics) cannot always be reproduced bit-for-bit through the
# property setter. The values chosen here are the inputs that land closest
# to (and in most cases exactly on) the recorded target widths once Excel
# re-quantizes them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List of items")

# Widen columns A:G to fit their (now longer) contents.
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668   # A -> 21
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666   # B -> ~14.43
$ws.Columns.Item(3).ColumnWidth = 13.166666666666666   # C -> 14
$ws.Columns.Item(4).ColumnWidth = 19                   # D -> ~19.86
$ws.Columns.Item(5).ColumnWidth = 5.666666666666667    # E -> ~6.43
$ws.Columns.Item(6).ColumnWidth = 8                    # F -> ~8.86
$ws.Columns.Item(7).ColumnWidth = 36.5                 # G -> ~37.29

# Move the live selection to where the user left off.
$ws.Activate()
$ws.Range("H23").Select()
